# Add a "logo path" column (D) to the teams sheet, one row per team,
# derived from the org name already present in column B for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 58

for ($r = 2; $r -le $lastRow; $r++) {
    $org = $ws.Cells.Item($r, 2).Value2
    if ($org -ne $null -and $org -ne "") {
        $ws.Cells.Item($r, 4).Value = "/static/logos/" + $org + ".png"
    }
}

# Restore the active selection to the cell the author left selected.
$ws.Range("I9").Select()
